$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il4"
$ws.Cells.Item(2, 3).Value = "Cd53"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.7922663333333334
$ws.Cells.Item(2, 8).Value = 2.376799
$ws.Cells.Item(2, 9).Value = 0.1759587713796512
$ws.Cells.Item(2, 10).Value = 0.1759587713796512
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.328017
$ws.Cells.Item(2, 14).Value = 0.984051
$ws.Cells.Item(2, 15).Value = 0.001744649556355686
$ws.Cells.Item(2, 16).Value = 0.001744649556355686
$ws.Cells.Item(2, 17).Value = 0.259876825861
$ws.Cells.Item(2, 18).Value = 2.338891432749
$ws.Cells.Item(2, 19).Value = 0.0003069863924244001
$ws.Cells.Item(2, 20).Value = 0.0003069863924244001

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il4"
$ws.Cells.Item(3, 3).Value = "Cd53"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.7922663333333334
$ws.Cells.Item(3, 8).Value = 2.376799
$ws.Cells.Item(3, 9).Value = 0.1759587713796512
$ws.Cells.Item(3, 10).Value = 0.1759587713796512
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 185.1323926666666
$ws.Cells.Item(3, 14).Value = 555.3971779999999
$ws.Cells.Item(3, 15).Value = 0.9846780707492804
$ws.Cells.Item(3, 16).Value = 0.9846780707492802
$ws.Cells.Item(3, 17).Value = 146.6741619192469
$ws.Cells.Item(3, 18).Value = 1320.067457273222
$ws.Cells.Item(3, 19).Value = 0.1732627435335286
$ws.Cells.Item(3, 20).Value = 0.1732627435335286

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il4"
$ws.Cells.Item(4, 3).Value = "Cd53"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.7922663333333334
$ws.Cells.Item(4, 8).Value = 2.376799
$ws.Cells.Item(4, 9).Value = 0.1759587713796512
$ws.Cells.Item(4, 10).Value = 0.1759587713796512
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 2.552706666666667
$ws.Cells.Item(4, 14).Value = 7.65812
$ws.Cells.Item(4, 15).Value = 0.01357727969436402
$ws.Cells.Item(4, 16).Value = 0.01357727969436402
$ws.Cells.Item(4, 17).Value = 2.022423550875556
$ws.Cells.Item(4, 18).Value = 18.20181195788
$ws.Cells.Item(4, 19).Value = 0.002389041453698179
$ws.Cells.Item(4, 20).Value = 0.002389041453698179

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Il4"
$ws.Cells.Item(5, 3).Value = "Cd53"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.077831666666667
$ws.Cells.Item(5, 8).Value = 6.233495
$ws.Cells.Item(5, 9).Value = 0.4614770208171574
$ws.Cells.Item(5, 10).Value = 0.4614770208171574
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.328017
$ws.Cells.Item(5, 14).Value = 0.984051
$ws.Cells.Item(5, 15).Value = 0.001744649556355686
$ws.Cells.Item(5, 16).Value = 0.001744649556355686
$ws.Cells.Item(5, 17).Value = 0.6815641098050002
$ws.Cells.Item(5, 18).Value = 6.134076988245001
$ws.Cells.Item(5, 19).Value = 0.0008051156796369974
$ws.Cells.Item(5, 20).Value = 0.0008051156796369974

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Il4"
$ws.Cells.Item(6, 3).Value = "Cd53"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.077831666666667
$ws.Cells.Item(6, 8).Value = 6.233495
$ws.Cells.Item(6, 9).Value = 0.4614770208171574
$ws.Cells.Item(6, 10).Value = 0.4614770208171574
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 185.1323926666666
$ws.Cells.Item(6, 14).Value = 555.3971779999999
$ws.Cells.Item(6, 15).Value = 0.9846780707492804
$ws.Cells.Item(6, 16).Value = 0.9846780707492802
$ws.Cells.Item(6, 17).Value = 384.6739480085678
$ws.Cells.Item(6, 18).Value = 3462.06553207711
$ws.Cells.Item(6, 19).Value = 0.454406302553364
$ws.Cells.Item(6, 20).Value = 0.4544063025533639

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Il4"
$ws.Cells.Item(7, 3).Value = "Cd53"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.077831666666667
$ws.Cells.Item(7, 8).Value = 6.233495
$ws.Cells.Item(7, 9).Value = 0.4614770208171574
$ws.Cells.Item(7, 10).Value = 0.4614770208171574
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 2.552706666666667
$ws.Cells.Item(7, 14).Value = 7.65812
$ws.Cells.Item(7, 15).Value = 0.01357727969436402
$ws.Cells.Item(7, 16).Value = 0.01357727969436402
$ws.Cells.Item(7, 17).Value = 5.304094747711113
$ws.Cells.Item(7, 18).Value = 47.73685272940001
$ws.Cells.Item(7, 19).Value = 0.006265602584156394
$ws.Cells.Item(7, 20).Value = 0.006265602584156393

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Il4"
$ws.Cells.Item(8, 3).Value = "Cd53"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.8480786666666668
$ws.Cells.Item(8, 8).Value = 2.544236
$ws.Cells.Item(8, 9).Value = 0.1883544383264543
$ws.Cells.Item(8, 10).Value = 0.1883544383264543
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.328017
$ws.Cells.Item(8, 14).Value = 0.984051
$ws.Cells.Item(8, 15).Value = 0.001744649556355686
$ws.Cells.Item(8, 16).Value = 0.001744649556355686
$ws.Cells.Item(8, 17).Value = 0.2781842200040001
$ws.Cells.Item(8, 18).Value = 2.503657980036
$ws.Cells.Item(8, 19).Value = 0.0003286124872638729
$ws.Cells.Item(8, 20).Value = 0.000328612487263873

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Il4"
$ws.Cells.Item(9, 3).Value = "Cd53"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.8480786666666668
$ws.Cells.Item(9, 8).Value = 2.544236
$ws.Cells.Item(9, 9).Value = 0.1883544383264543
$ws.Cells.Item(9, 10).Value = 0.1883544383264543
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 185.1323926666666
$ws.Cells.Item(9, 14).Value = 555.3971779999999
$ws.Cells.Item(9, 15).Value = 0.9846780707492804
$ws.Cells.Item(9, 16).Value = 0.9846780707492802
$ws.Cells.Item(9, 17).Value = 157.0068327295564
$ws.Cells.Item(9, 18).Value = 1413.061494566008
$ws.Cells.Item(9, 19).Value = 0.1854684849483573
$ws.Cells.Item(9, 20).Value = 0.1854684849483573

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Il4"
$ws.Cells.Item(10, 3).Value = "Cd53"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8480786666666668
$ws.Cells.Item(10, 8).Value = 2.544236
$ws.Cells.Item(10, 9).Value = 0.1883544383264543
$ws.Cells.Item(10, 10).Value = 0.1883544383264543
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 2.552706666666667
$ws.Cells.Item(10, 14).Value = 7.65812
$ws.Cells.Item(10, 15).Value = 0.01357727969436402
$ws.Cells.Item(10, 16).Value = 0.01357727969436402
$ws.Cells.Item(10, 17).Value = 2.164896066257778
$ws.Cells.Item(10, 18).Value = 19.48406459632
$ws.Cells.Item(10, 19).Value = 0.002557340890833108
$ws.Cells.Item(10, 20).Value = 0.002557340890833108

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Il4"
$ws.Cells.Item(11, 3).Value = "Cd53"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.7843913333333336
$ws.Cells.Item(11, 8).Value = 2.353174000000001
$ws.Cells.Item(11, 9).Value = 0.1742097694767371
$ws.Cells.Item(11, 10).Value = 0.1742097694767372
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.328017
$ws.Cells.Item(11, 14).Value = 0.984051
$ws.Cells.Item(11, 15).Value = 0.001744649556355686
$ws.Cells.Item(11, 16).Value = 0.001744649556355686
$ws.Cells.Item(11, 17).Value = 0.2572936919860001
$ws.Cells.Item(11, 18).Value = 2.315643227874001
$ws.Cells.Item(11, 19).Value = 0.0003039349970304159
$ws.Cells.Item(11, 20).Value = 0.0003039349970304159

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Il4"
$ws.Cells.Item(12, 3).Value = "Cd53"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.7843913333333336
$ws.Cells.Item(12, 8).Value = 2.353174000000001
$ws.Cells.Item(12, 9).Value = 0.1742097694767371
$ws.Cells.Item(12, 10).Value = 0.1742097694767372
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 185.1323926666666
$ws.Cells.Item(12, 14).Value = 555.3971779999999
$ws.Cells.Item(12, 15).Value = 0.9846780707492804
$ws.Cells.Item(12, 16).Value = 0.9846780707492802
$ws.Cells.Item(12, 17).Value = 145.2162443269969
$ws.Cells.Item(12, 18).Value = 1306.946198942972
$ws.Cells.Item(12, 19).Value = 0.1715405397140304
$ws.Cells.Item(12, 20).Value = 0.1715405397140304

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Il4"
$ws.Cells.Item(13, 3).Value = "Cd53"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.7843913333333336
$ws.Cells.Item(13, 8).Value = 2.353174000000001
$ws.Cells.Item(13, 9).Value = 0.1742097694767371
$ws.Cells.Item(13, 10).Value = 0.1742097694767372
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 2.552706666666667
$ws.Cells.Item(13, 14).Value = 7.65812
$ws.Cells.Item(13, 15).Value = 0.01357727969436402
$ws.Cells.Item(13, 16).Value = 0.01357727969436402
$ws.Cells.Item(13, 17).Value = 2.002320985875556
$ws.Cells.Item(13, 18).Value = 18.02088887288
$ws.Cells.Item(13, 19).Value = 0.00236529476567634
$ws.Cells.Item(13, 20).Value = 0.00236529476567634
